$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header row values for the new Wins/Losses/Ties columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style (bold, centered, bordered) from an existing header cell (AC1)
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill in the team record values for every data row (2 through 49)
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 77
    $ws.Cells.Item($r, 31).Value = 85
    $ws.Cells.Item($r, 32).Value = 0
}
